$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create new row 121 by copying formatting from row 120, then fill in values ---
$ws.Range("A120:V120").Copy()
$ws.Range("A121:V121").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 79: columns F-V (A-E unchanged)
$ws.Range("F79").Value = "Kasimpasa"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = "Fenerbahce"
$ws.Range("I79").Value = 2
$ws.Range("J79").Value = 4.89
$ws.Range("K79").Value = "01/10/2023 17:13"
$ws.Range("L79").Value = 6.87
$ws.Range("M79").Value = "08/10/2023 17:59"
$ws.Range("N79").Value = 4.43
$ws.Range("O79").Value = "01/10/2023 17:13"
$ws.Range("P79").Value = 4.9
$ws.Range("Q79").Value = "08/10/2023 17:59"
$ws.Range("R79").Value = 1.63
$ws.Range("S79").Value = "01/10/2023 17:13"
$ws.Range("T79").Value = 1.47
$ws.Range("U79").Value = "08/10/2023 17:59"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/turkey/super-lig/kasimpasa-fenerbahce/WQZd0xgp/"

# Row 80: columns F-V (A-E unchanged)
$ws.Range("F80").Value = "Hatayspor"
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = "Konyaspor"
$ws.Range("I80").Value = 1
$ws.Range("J80").Value = 2.54
$ws.Range("K80").Value = "01/10/2023 17:13"
$ws.Range("L80").Value = 2.54
$ws.Range("M80").Value = "08/10/2023 17:56"
$ws.Range("N80").Value = 3.46
$ws.Range("O80").Value = "01/10/2023 17:13"
$ws.Range("P80").Value = 3.16
$ws.Range("Q80").Value = "08/10/2023 17:57"
$ws.Range("R80").Value = 2.86
$ws.Range("S80").Value = "01/10/2023 17:13"
$ws.Range("T80").Value = 3.15
$ws.Range("U80").Value = "08/10/2023 17:56"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-konyaspor/0rqR5bVM/"

# Row 97: columns F-V (A-E unchanged)
$ws.Range("F97").Value = "Pendikspor"
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = "Fenerbahce"
$ws.Range("I97").Value = 5
$ws.Range("J97").Value = 7.87
$ws.Range("K97").Value = "23/10/2023 05:42"
$ws.Range("L97").Value = 10.56
$ws.Range("M97").Value = "29/10/2023 16:59"
$ws.Range("N97").Value = 5.68
$ws.Range("O97").Value = "23/10/2023 05:42"
$ws.Range("P97").Value = 6.22
$ws.Range("Q97").Value = "29/10/2023 16:59"
$ws.Range("R97").Value = 1.35
$ws.Range("S97").Value = "23/10/2023 05:42"
$ws.Range("T97").Value = 1.28
$ws.Range("U97").Value = "29/10/2023 16:59"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/turkey/super-lig/pendikspor-fenerbahce/vc8IQY6k/"

# Row 98: columns F-V (A-E unchanged)
$ws.Range("F98").Value = "Antalyaspor"
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = "Basaksehir"
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2.03
$ws.Range("K98").Value = "22/10/2023 20:15"
$ws.Range("L98").Value = 2.11
$ws.Range("M98").Value = "29/10/2023 16:54"
$ws.Range("N98").Value = 3.51
$ws.Range("O98").Value = "22/10/2023 20:15"
$ws.Range("P98").Value = 3.35
$ws.Range("Q98").Value = "29/10/2023 16:54"
$ws.Range("R98").Value = 3.8
$ws.Range("S98").Value = "22/10/2023 20:15"
$ws.Range("T98").Value = 3.89
$ws.Range("U98").Value = "29/10/2023 16:54"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-basaksehir/tz8GozqF/"

# Row 99: columns F-V (A-E unchanged)
$ws.Range("F99").Value = "Ankaragucu"
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = "Samsunspor"
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2.19
$ws.Range("K99").Value = "23/10/2023 05:42"
$ws.Range("L99").Value = 2.7
$ws.Range("M99").Value = "29/10/2023 16:59"
$ws.Range("N99").Value = 3.59
$ws.Range("O99").Value = "23/10/2023 05:42"
$ws.Range("P99").Value = 3.32
$ws.Range("Q99").Value = "29/10/2023 16:54"
$ws.Range("R99").Value = 3.29
$ws.Range("S99").Value = "23/10/2023 05:42"
$ws.Range("T99").Value = 2.82
$ws.Range("U99").Value = "29/10/2023 16:59"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/turkey/super-lig/ankaragucu-samsunspor/2kdPqEDR/"

# Row 115: columns F-V (A-E unchanged)
$ws.Range("F115").Value = "Alanyaspor"
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = "Gaziantep"
$ws.Range("I115").Value = 3
$ws.Range("J115").Value = 2.07
$ws.Range("K115").Value = "06/11/2023 18:12"
$ws.Range("L115").Value = 2.34
$ws.Range("M115").Value = "11/11/2023 13:59"
$ws.Range("N115").Value = 3.62
$ws.Range("O115").Value = "06/11/2023 18:12"
$ws.Range("P115").Value = 3.29
$ws.Range("Q115").Value = "11/11/2023 13:55"
$ws.Range("R115").Value = 3.67
$ws.Range("S115").Value = "06/11/2023 18:12"
$ws.Range("T115").Value = 3.37
$ws.Range("U115").Value = "11/11/2023 13:59"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-gaziantep/hb1qysRO/"

# Row 116: columns F-V (A-E unchanged)
$ws.Range("F116").Value = "Karagumruk"
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = "Sivasspor"
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2.41
$ws.Range("K116").Value = "06/11/2023 04:12"
$ws.Range("L116").Value = 2.22
$ws.Range("M116").Value = "11/11/2023 13:53"
$ws.Range("N116").Value = 3.41
$ws.Range("O116").Value = "06/11/2023 04:12"
$ws.Range("P116").Value = 3.42
$ws.Range("Q116").Value = "11/11/2023 13:53"
$ws.Range("R116").Value = 3.1
$ws.Range("S116").Value = "06/11/2023 04:12"
$ws.Range("T116").Value = 3.51
$ws.Range("U116").Value = "11/11/2023 13:53"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/turkey/super-lig/f-karagumruk-sivasspor/WOWgepCt/"

# Row 119: columns F-V (A-E unchanged)
$ws.Range("F119").Value = "Rizespor"
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = "Istanbulspor AS"
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1.8
$ws.Range("K119").Value = "06/11/2023 18:12"
$ws.Range("L119").Value = 1.69
$ws.Range("M119").Value = "12/11/2023 13:53"
$ws.Range("N119").Value = 3.94
$ws.Range("O119").Value = "06/11/2023 18:12"
$ws.Range("P119").Value = 4.01
$ws.Range("Q119").Value = "12/11/2023 13:53"
$ws.Range("R119").Value = 4.32
$ws.Range("S119").Value = "06/11/2023 18:12"
$ws.Range("T119").Value = 5.3
$ws.Range("U119").Value = "12/11/2023 13:52"
$ws.Range("V119").Value = "https://www.betexplorer.com/football/turkey/super-lig/rizespor-istanbulspor-as/04OAi2B5/"

# Row 120: columns F-V (A-E unchanged)
$ws.Range("F120").Value = "Besiktas"
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = "Basaksehir"
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 1.55
$ws.Range("K120").Value = "05/11/2023 17:12"
$ws.Range("L120").Value = 1.74
$ws.Range("M120").Value = "12/11/2023 13:56"
$ws.Range("N120").Value = 4.44
$ws.Range("O120").Value = "05/11/2023 17:12"
$ws.Range("P120").Value = 3.84
$ws.Range("Q120").Value = "12/11/2023 13:59"
$ws.Range("R120").Value = 5.79
$ws.Range("S120").Value = "05/11/2023 17:12"
$ws.Range("T120").Value = 5.16
$ws.Range("U120").Value = "12/11/2023 13:59"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/turkey/super-lig/besiktas-basaksehir/vkz2gOtg/"

# Row 121: new row, set all columns A-V
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "turkey"
$ws.Range("C121").Value = "super-lig"
$ws.Range("D121").Value = "2023-2024"
$ws.Range("E121").Value = 45242.70833333334
$ws.Range("F121").Value = "Adana Demirspor"
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = "Fenerbahce"
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3.48
$ws.Range("K121").Value = "05/11/2023 17:12"
$ws.Range("L121").Value = 3.54
$ws.Range("M121").Value = "12/11/2023 16:57"
$ws.Range("N121").Value = 4.1
$ws.Range("O121").Value = "05/11/2023 17:12"
$ws.Range("P121").Value = 4.42
$ws.Range("Q121").Value = "12/11/2023 16:57"
$ws.Range("R121").Value = 2
$ws.Range("S121").Value = "05/11/2023 17:12"
$ws.Range("T121").Value = 1.92
$ws.Range("U121").Value = "12/11/2023 16:57"
$ws.Range("V121").Value = "https://www.betexplorer.com/football/turkey/super-lig/adanademirspor-fenerbahce/UFpSaQ3P/"
